# Add two new columns, I (I0) and J (IF), mirroring the existing
# header style used by column H (IP).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - set values first, then copy the formatting from the
# existing H1 header cell (xlPasteFormats = -4122) so the new headers
# reuse the same bold/bordered/centered style as the rest of row 1
# instead of minting a new style entry.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 4

$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
